# Scheduled runner update: refresh Leve market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) pulled from the
# Universalis API snapshot, across the ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value2 = 931.6667
$ws.Cells.Item(70, 9).Value2 = 1000
$ws.Cells.Item(70, 10).Value2 = 918
$ws.Cells.Item(70, 11).Value2 = 3000
$ws.Cells.Item(70, 12).Value2 = 2754
$ws.Cells.Item(70, 13).Value2 = -2730
$ws.Cells.Item(70, 14).Value2 = -3294
$ws.Cells.Item(73, 8).Value2 = 931.6667
$ws.Cells.Item(73, 9).Value2 = 1000
$ws.Cells.Item(73, 10).Value2 = 918
$ws.Cells.Item(73, 11).Value2 = 3000
$ws.Cells.Item(73, 12).Value2 = 2754
$ws.Cells.Item(73, 13).Value2 = -2064
$ws.Cells.Item(73, 14).Value2 = -4626
$ws.Cells.Item(88, 8).Value2 = 1933.1818
$ws.Cells.Item(88, 10).Value2 = 1914.375
$ws.Cells.Item(88, 12).Value2 = 1914.375
$ws.Cells.Item(88, 14).Value2 = -2726.375
$ws.Cells.Item(91, 8).Value2 = 1933.1818
$ws.Cells.Item(91, 10).Value2 = 1914.375
$ws.Cells.Item(91, 12).Value2 = 1914.375
$ws.Cells.Item(91, 14).Value2 = -4722.375
$ws.Cells.Item(94, 8).Value2 = 1749.75
$ws.Cells.Item(94, 9).Value2 = 1749.75
$ws.Cells.Item(94, 10).Value2 = 0
$ws.Cells.Item(94, 11).Value2 = 1749.75
$ws.Cells.Item(94, 12).Value2 = 0
$ws.Cells.Item(94, 13).Value2 = -1298.75
$ws.Cells.Item(94, 14).ClearContents()
$ws.Cells.Item(97, 8).Value2 = 1068.2
$ws.Cells.Item(97, 10).Value2 = 1068.2
$ws.Cells.Item(97, 12).Value2 = 3204.6
$ws.Cells.Item(97, 14).Value2 = -4196.6
$ws.Cells.Item(107, 8).Value2 = 378.21738
$ws.Cells.Item(107, 9).Value2 = 358.07693
$ws.Cells.Item(107, 10).Value2 = 404.4
$ws.Cells.Item(107, 11).Value2 = 358.07693
$ws.Cells.Item(107, 12).Value2 = 404.4
$ws.Cells.Item(107, 13).Value2 = 1561.92307
$ws.Cells.Item(107, 14).Value2 = -4244.4
$ws.Cells.Item(132, 8).Value2 = 1939.3334
$ws.Cells.Item(132, 9).Value2 = 1938.0426
$ws.Cells.Item(132, 11).Value2 = 5814.1278
$ws.Cells.Item(132, 13).Value2 = -3284.1278
$ws.Cells.Item(137, 8).Value2 = 2660.0417
$ws.Cells.Item(137, 9).Value2 = 2477.05
$ws.Cells.Item(137, 11).Value2 = 7431.150000000001
$ws.Cells.Item(137, 13).Value2 = -4881.150000000001
$ws.Cells.Item(138, 8).Value2 = 1853.2222
$ws.Cells.Item(138, 9).Value2 = 769.0454999999999
$ws.Cells.Item(138, 10).Value2 = 2257.4915
$ws.Cells.Item(138, 11).Value2 = 2307.1365
$ws.Cells.Item(138, 12).Value2 = 6772.4745
$ws.Cells.Item(138, 13).Value2 = 2832.8635
$ws.Cells.Item(138, 14).Value2 = -17052.4745

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value2 = 1671.8667
$ws.Cells.Item(122, 9).Value2 = 1320.1111
$ws.Cells.Item(122, 10).Value2 = 2199.5
$ws.Cells.Item(122, 11).Value2 = 3960.3333
$ws.Cells.Item(122, 12).Value2 = 6598.5
$ws.Cells.Item(122, 13).Value2 = -1510.3333
$ws.Cells.Item(122, 14).Value2 = -11498.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value2 = 1148.5714
$ws.Cells.Item(16, 9).Value2 = 1208
$ws.Cells.Item(16, 10).Value2 = 1000
$ws.Cells.Item(16, 11).Value2 = 1208
$ws.Cells.Item(16, 12).Value2 = 1000
$ws.Cells.Item(16, 13).Value2 = -921
$ws.Cells.Item(16, 14).Value2 = -1574
$ws.Cells.Item(31, 8).Value2 = 11119.021
$ws.Cells.Item(31, 9).Value2 = 15542.857
$ws.Cells.Item(31, 10).Value2 = 4599.684
$ws.Cells.Item(31, 11).Value2 = 15542.857
$ws.Cells.Item(31, 12).Value2 = 4599.684
$ws.Cells.Item(31, 13).Value2 = -15247.857
$ws.Cells.Item(31, 14).Value2 = -5189.684
$ws.Cells.Item(34, 8).Value2 = 11119.021
$ws.Cells.Item(34, 9).Value2 = 15542.857
$ws.Cells.Item(34, 10).Value2 = 4599.684
$ws.Cells.Item(34, 11).Value2 = 15542.857
$ws.Cells.Item(34, 12).Value2 = 4599.684
$ws.Cells.Item(34, 13).Value2 = -15340.857
$ws.Cells.Item(34, 14).Value2 = -5003.684
$ws.Cells.Item(113, 8).Value2 = 1148.5714
$ws.Cells.Item(113, 9).Value2 = 1208
$ws.Cells.Item(113, 10).Value2 = 1000
$ws.Cells.Item(113, 11).Value2 = 1208
$ws.Cells.Item(113, 12).Value2 = 1000
$ws.Cells.Item(113, 13).Value2 = 962
$ws.Cells.Item(113, 14).Value2 = -5340
$ws.Cells.Item(134, 8).Value2 = 1265.2982
$ws.Cells.Item(134, 9).Value2 = 894.35297
$ws.Cells.Item(134, 10).Value2 = 1813.6522
$ws.Cells.Item(134, 11).Value2 = 2683.05891
$ws.Cells.Item(134, 12).Value2 = 5440.9566
$ws.Cells.Item(134, 13).Value2 = -148.0589100000002
$ws.Cells.Item(134, 14).Value2 = -10510.9566

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value2 = 6033.1113
$ws.Cells.Item(62, 9).Value2 = 1872
$ws.Cells.Item(62, 10).Value2 = 9362
$ws.Cells.Item(62, 11).Value2 = 5616
$ws.Cells.Item(62, 12).Value2 = 28086
$ws.Cells.Item(62, 13).Value2 = -4930
$ws.Cells.Item(62, 14).Value2 = -29458
$ws.Cells.Item(65, 8).Value2 = 6033.1113
$ws.Cells.Item(65, 9).Value2 = 1872
$ws.Cells.Item(65, 10).Value2 = 9362
$ws.Cells.Item(65, 11).Value2 = 16848
$ws.Cells.Item(65, 12).Value2 = 84258
$ws.Cells.Item(65, 13).Value2 = -13416
$ws.Cells.Item(65, 14).Value2 = -91122
$ws.Cells.Item(131, 8).Value2 = 778.16
$ws.Cells.Item(131, 10).Value2 = 778.16
$ws.Cells.Item(131, 12).Value2 = 2334.48
$ws.Cells.Item(131, 14).Value2 = -12414.48
$ws.Cells.Item(132, 8).Value2 = 1638.75
$ws.Cells.Item(132, 9).Value2 = 555
$ws.Cells.Item(132, 10).Value2 = 2000
$ws.Cells.Item(132, 11).Value2 = 4995
$ws.Cells.Item(132, 12).Value2 = 18000
$ws.Cells.Item(132, 13).Value2 = -2465
$ws.Cells.Item(132, 14).Value2 = -23060
$ws.Cells.Item(136, 8).Value2 = 828
$ws.Cells.Item(136, 9).Value2 = 828
$ws.Cells.Item(136, 11).Value2 = 2484
$ws.Cells.Item(136, 13).Value2 = 2616

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value2 = 44445464
$ws.Cells.Item(122, 9).Value2 = 16667782
$ws.Cells.Item(122, 11).Value2 = 50003346
$ws.Cells.Item(122, 13).Value2 = -50000896
$ws.Cells.Item(132, 8).Value2 = 20975.793
$ws.Cells.Item(132, 9).Value2 = 3847.5
$ws.Cells.Item(132, 10).Value2 = 74807.57000000001
$ws.Cells.Item(132, 11).Value2 = 11542.5
$ws.Cells.Item(132, 12).Value2 = 224422.71
$ws.Cells.Item(132, 13).Value2 = -9012.5
$ws.Cells.Item(132, 14).Value2 = -229482.71

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value2 = 2555.7334
$ws.Cells.Item(22, 10).Value2 = 3866.3333
$ws.Cells.Item(22, 12).Value2 = 3866.3333
$ws.Cells.Item(22, 14).Value2 = -4456.3333
$ws.Cells.Item(27, 8).Value2 = 2555.7334
$ws.Cells.Item(27, 10).Value2 = 3866.3333
$ws.Cells.Item(27, 12).Value2 = 3866.3333
$ws.Cells.Item(27, 14).Value2 = -4080.3333
$ws.Cells.Item(46, 8).Value2 = 773.3333
$ws.Cells.Item(46, 9).Value2 = 698.7143
$ws.Cells.Item(46, 11).Value2 = 698.7143
$ws.Cells.Item(46, 13).Value2 = -510.7143
$ws.Cells.Item(64, 8).Value2 = 30000
$ws.Cells.Item(64, 10).Value2 = 30000
$ws.Cells.Item(64, 12).Value2 = 30000
$ws.Cells.Item(64, 14).Value2 = -30450
$ws.Cells.Item(67, 8).Value2 = 30000
$ws.Cells.Item(67, 10).Value2 = 30000
$ws.Cells.Item(67, 12).Value2 = 30000
$ws.Cells.Item(67, 14).Value2 = -31560
$ws.Cells.Item(93, 8).Value2 = 1755.0476
$ws.Cells.Item(93, 9).Value2 = 1810.3529
$ws.Cells.Item(93, 10).Value2 = 1520
$ws.Cells.Item(93, 11).Value2 = 1810.3529
$ws.Cells.Item(93, 12).Value2 = 1520
$ws.Cells.Item(93, 13).Value2 = -562.3529000000001
$ws.Cells.Item(93, 14).Value2 = -4016
$ws.Cells.Item(132, 8).Value2 = 2553.4375
$ws.Cells.Item(132, 9).Value2 = 1233.7222
$ws.Cells.Item(132, 11).Value2 = 3701.1666
$ws.Cells.Item(132, 13).Value2 = -1171.1666
$ws.Cells.Item(136, 8).Value2 = 23708
$ws.Cells.Item(136, 9).Value2 = 32241.938
$ws.Cells.Item(136, 10).Value2 = 950.8333
$ws.Cells.Item(136, 11).Value2 = 96725.814
$ws.Cells.Item(136, 12).Value2 = 2852.4999
$ws.Cells.Item(136, 13).Value2 = -94175.814
$ws.Cells.Item(136, 14).Value2 = -7952.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value2 = 0
$ws.Cells.Item(17, 9).Value2 = 0
$ws.Cells.Item(17, 11).Value2 = 0
$ws.Cells.Item(17, 13).ClearContents()
$ws.Cells.Item(62, 8).Value2 = 4752.7334
$ws.Cells.Item(62, 9).Value2 = 3964.6
$ws.Cells.Item(62, 10).Value2 = 5146.8
$ws.Cells.Item(62, 11).Value2 = 3964.6
$ws.Cells.Item(62, 12).Value2 = 5146.8
$ws.Cells.Item(62, 13).Value2 = -3340.6
$ws.Cells.Item(62, 14).Value2 = -6394.8
$ws.Cells.Item(63, 8).Value2 = 39997.5
$ws.Cells.Item(63, 10).Value2 = 39997.5
$ws.Cells.Item(63, 12).Value2 = 39997.5
$ws.Cells.Item(63, 14).Value2 = -41245.5
$ws.Cells.Item(65, 8).Value2 = 4752.7334
$ws.Cells.Item(65, 9).Value2 = 3964.6
$ws.Cells.Item(65, 10).Value2 = 5146.8
$ws.Cells.Item(65, 11).Value2 = 19823
$ws.Cells.Item(65, 12).Value2 = 25734
$ws.Cells.Item(65, 13).Value2 = -16703
$ws.Cells.Item(65, 14).Value2 = -31974
$ws.Cells.Item(66, 8).Value2 = 39997.5
$ws.Cells.Item(66, 10).Value2 = 39997.5
$ws.Cells.Item(66, 12).Value2 = 119992.5
$ws.Cells.Item(66, 14).Value2 = -126232.5
$ws.Cells.Item(96, 8).Value2 = 1440
$ws.Cells.Item(96, 9).Value2 = 1425
$ws.Cells.Item(96, 10).Value2 = 1500
$ws.Cells.Item(96, 11).Value2 = 1425
$ws.Cells.Item(96, 12).Value2 = 1500
$ws.Cells.Item(96, 13).Value2 = -52
$ws.Cells.Item(96, 14).Value2 = -4246
$ws.Cells.Item(136, 8).Value2 = 38463740
$ws.Cells.Item(136, 9).Value2 = 71430730
$ws.Cells.Item(136, 11).Value2 = 214292190
$ws.Cells.Item(136, 13).Value2 = -214289640
